$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3 all shared the same string)
# ---------------------------------------------------------------------------
$handedBack = "Handed back: in sync with en-US"

$wsOverview.Range("E2").Value = $handedBack
$wsOverview.Range("F2").Value = $handedBack
$wsOverview.Range("E3").Value = $handedBack
$wsOverview.Range("F3").Value = $handedBack

$wsZh.Range("C2").Value = $handedBack
$wsZh.Range("C3").Value = $handedBack

$wsDe.Range("C2").Value = $handedBack
$wsDe.Range("C3").Value = $handedBack

# ---------------------------------------------------------------------------
# 2. Populate the "Latest Target File" (I) and "Latest Handback File" (J)
#    columns on the zh-cn and de-de report sheets, plus the de-de
#    "Latest Handback DateTime" (K) column - the handback for de-de
#    completed, zh-cn has not been handed back yet.
# ---------------------------------------------------------------------------

# -- zh-cn --
$wsZh.Range("I2").Value = "1daa9765-d5e8-4f9f-a8cd-e589105224c0.md"
$wsZh.Range("J2").Value = "1daa9765-d5e8-4f9f-a8cd-e589105224c0.dc57d2b065390c90d27816386da34daf2d8b263f.zh-cn.xlf"

$wsZh.Range("I3").Value = "a76b8f79-3bfc-467c-80fd-cf3c63741b28.md"
$wsZh.Range("J3").Value = "a76b8f79-3bfc-467c-80fd-cf3c63741b28.31bd3791c66b9e7ee668cc36dde2633eab440f78.zh-cn.xlf"

# -- de-de --
$wsDe.Range("I2").Value = "1daa9765-d5e8-4f9f-a8cd-e589105224c0.md"
$wsDe.Range("J2").Value = "1daa9765-d5e8-4f9f-a8cd-e589105224c0.dc57d2b065390c90d27816386da34daf2d8b263f.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-29 10:48:55"

$wsDe.Range("I3").Value = "a76b8f79-3bfc-467c-80fd-cf3c63741b28.md"
$wsDe.Range("J3").Value = "a76b8f79-3bfc-467c-80fd-cf3c63741b28.31bd3791c66b9e7ee668cc36dde2633eab440f78.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-29 10:48:55"

# ---------------------------------------------------------------------------
# 3. Hyperlink style for the new "Latest Target File" cells, matching the
#    existing "Source File Name" hyperlink look (cell style 1).
# ---------------------------------------------------------------------------
$wsZh.Range("I2").Style = $wsZh.Range("A2").Style
$wsZh.Range("I3").Style = $wsZh.Range("A3").Style
$wsDe.Range("I2").Style = $wsDe.Range("A2").Style
$wsDe.Range("I3").Style = $wsDe.Range("A3").Style

# ---------------------------------------------------------------------------
# 4. Hyperlinks - rebuild so the "Latest Target File" cells (I2/I3) link to
#    the same commit-pinned GitHub blob URLs as the "Source File Name"
#    cells (A2/A3), preserving the existing link targets/order.
# ---------------------------------------------------------------------------
$addr1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f695303534129e140d551b0a0e6aa92cbc97e0ab/e2e/1daa9765-d5e8-4f9f-a8cd-e589105224c0.md"
$disp1 = "1daa9765-d5e8-4f9f-a8cd-e589105224c0.md"
$addr2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f695303534129e140d551b0a0e6aa92cbc97e0ab/e2e/a76b8f79-3bfc-467c-80fd-cf3c63741b28.md"
$disp2 = "a76b8f79-3bfc-467c-80fd-cf3c63741b28.md"

function Rebuild-Hyperlinks($ws) {
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $addr1, "", "", $disp1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $addr1, "", "", $disp1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $addr2, "", "", $disp2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $addr2, "", "", $disp2)
}

Rebuild-Hyperlinks $wsZh
Rebuild-Hyperlinks $wsDe

# ---------------------------------------------------------------------------
# 5. Column widths - the Status/Latest Target File/Latest Handback File
#    columns were widened to fit the new report content.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
